$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New row 11 - SRS_02
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "SRS_02"
$ws.Range("B11").Value = "Software Requirements Specification"
$ws.Range("C11").Value = "Nourhan Ali,Mohamed Ibrahim`n,Manar Ali ,Aalaa Adel and Al-Shimaa`nShehata"
$ws.Range("D11").Value = "Mostafa Mohamed"
$ws.Range("E11").Value = 45508
$ws.Range("F11").Value = "__"
$ws.Range("G11").Value = "__"
$ws.Range("H11").Value = "Solved"

# ---------------------------------------------------------------------
# New row 12 - RTM_01
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "RTM_01"
$ws.Range("B12").Value = "Requirment Tracability Matrix"
$ws.Range("C12").Value = "Mostafa Mohamed"
$ws.Range("D12").Value = "Mohamed Ibrahim"
$ws.Range("E12").Value = 45508
$ws.Range("F12").Value = "__"
$ws.Range("G12").Value = "__"
$ws.Range("H12").Value = "Solved"

# ---------------------------------------------------------------------
# Formatting - reuse the look of the existing data rows so the new rows
# match the rest of the table (center/middle aligned, wrapped where the
# existing "problems"/"recommended resolution" columns wrap, bold status
# column, etc.)
# ---------------------------------------------------------------------

# A11,B11,D11,G11 -> same look as B4 (center + middle, no wrap)
$ws.Range("B4").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null

# C11,F11 -> same look as C4 (center + middle + wrap)
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null

# H11,H12 -> same look as H4 (bold status font, center + middle)
$ws.Range("H4").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null

# A12 -> same look as A5 (center only)
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null

# B12,C12,D12,G12 -> same look as B5 (center + middle, no wrap)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Date columns - short date, builtin format
$ws.Range("E11").NumberFormat = "mm-dd-yy"
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("E11").VerticalAlignment = -4108

$ws.Range("E12").NumberFormat = "mm-dd-yy"
$ws.Range("E12").HorizontalAlignment = -4108

# F12 - centered, top aligned, wrapped
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("F12").VerticalAlignment = -4160
$ws.Range("F12").WrapText = $true

# ---------------------------------------------------------------------
# Row heights to match the taller wrapped content
# ---------------------------------------------------------------------
$ws.Rows("11:11").RowHeight = 72
$ws.Rows("12:12").RowHeight = 21

# ---------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 109
$ws.Range("A13").Select() | Out-Null
